$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2, $colStart, $colEnd) {
    for ($c = $colStart; $c -le $colEnd; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value()
        $v2 = $ws.Cells.Item($r2, $c).Value()
        $ws.Cells.Item($r1, $c).Value = $v2
        $ws.Cells.Item($r2, $c).Value = $v1
    }
}

# Swap the content (columns B..AC, i.e. 2..29) of the following row pairs,
# keeping column A (the running index) untouched.
Swap-Rows 99 100 2 29
Swap-Rows 111 112 2 29
Swap-Rows 122 123 2 29

# Helper to set a full data row. $vals is a hashtable keyed by
# column letter (A..AC). Missing keys are simply left blank.
function Set-RowData($r, $vals) {
    $styleSrcA = $ws.Range("A133")
    $styleSrcE = $ws.Range("E133")

    $styleSrcA.Copy()
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null
    $styleSrcE.Copy()
    $ws.Range("E$r").PasteSpecial(-4122) | Out-Null

    foreach ($col in $vals.Keys) {
        $ws.Range("$col$r").Value = $vals[$col]
    }
}

$row134 = [ordered]@{
    A = 132
    B = 7952732
    C = "Bosnia Herzegovina Premier Liga"
    D = "Bosnia  Herzegovina Premier Liga"
    E = 45367.39583333334
    F = "Sloga"
    G = "FK Tuzla City"
    H = 3
    I = 1
    J = "H"
    K = 1.909
    L = 3.2
    M = 3.6
    N = 2.05
    O = 3.2
    P = 3.2
    Q = -0.25
    R = 1.825
    S = 1.975
    T = 2.25
    U = 1.925
    V = 1.875
    W = 1.05
    X = -1
    Y = -1
    Z = 0.825
    AA = -1
    AB = 0.925
    AC = -1
}
Set-RowData 134 $row134

$row135 = [ordered]@{
    A = 133
    B = 7952455
    C = "Bosnia Herzegovina Premier Liga"
    D = "Bosnia  Herzegovina Premier Liga"
    E = 45367.5
    F = "Zrinjski Mostar"
    G = "Zeljeznicar"
    H = 2
    I = 0
    J = "H"
    K = 1.363
    L = 4.2
    M = 7
    N = 1.333
    O = 3.6
    P = 11
    Q = -1.25
    R = 1.9
    S = 1.9
    T = 2.25
    U = 1.95
    V = 1.85
    W = 0.333
    X = -1
    Y = -1
    Z = 0.8999999999999999
    AA = -1
    AB = -0.5
    AC = 0.425
}
Set-RowData 135 $row135

$row136 = [ordered]@{
    A = 134
    B = 7952730
    C = "Bosnia Herzegovina Premier Liga"
    D = "Bosnia  Herzegovina Premier Liga"
    E = 45367.69791666666
    F = "FK Sarajevo"
    G = "GOSK Gabela"
    H = 4
    I = 0
    J = "H"
    K = 1.4
    L = 4
    M = 6.5
    N = 1.3
    O = 4.5
    P = 8
    Q = -1.5
    R = 1.95
    S = 1.85
    T = 2.75
    U = 1.9
    V = 1.9
    W = 0.3
    X = -1
    Y = -1
    Z = 0.95
    AA = -1
    AB = 0.8999999999999999
    AC = -1
}
Set-RowData 136 $row136

$row137 = [ordered]@{
    A = 135
    B = 7952731
    C = "Bosnia Herzegovina Premier Liga"
    D = "Bosnia  Herzegovina Premier Liga"
    E = 45368.41666666666
    F = "NK Posusje"
    G = "Siroki Brijeg"
    K = 1.666
    L = 3.2
    M = 5
    N = 1.65
    O = 3.2
    P = 5.5
    Q = -0.75
    R = 1.9
    S = 1.9
    T = 1.75
    U = 1.775
    V = 2.025
    W = 0
    X = 0
    Y = 0
    Z = 0
    AA = 0
}
Set-RowData 137 $row137

$row138 = [ordered]@{
    A = 136
    B = 7952734
    C = "Bosnia Herzegovina Premier Liga"
    D = "Bosnia  Herzegovina Premier Liga"
    E = 45368.58333333334
    F = "Borac Banja Luka"
    G = "Zvijezda 09"
    K = 1.181
    L = 5.5
    M = 11
    N = 1.2
    O = 5
    P = 11
    Q = -1.75
    R = 1.925
    S = 1.875
    T = 2.75
    U = 1.9
    V = 1.9
    W = 0
    X = 0
    Y = 0
    Z = 0
    AA = 0
}
Set-RowData 138 $row138

$row139 = [ordered]@{
    A = 137
    B = 7952733
    C = "Bosnia Herzegovina Premier Liga"
    D = "Bosnia  Herzegovina Premier Liga"
    E = 45368.66666666666
    F = "Velez Mostar"
    G = "NK Igman Konjic"
    K = 1.3
    L = 4.333
    M = 8.5
    N = 1.3
    O = 4.2
    P = 9.5
    Q = -1.25
    R = 1.825
    S = 1.975
    T = 2.25
    U = 1.775
    V = 2.025
    W = 0
    X = 0
    Y = 0
    Z = 0
    AA = 0
}
Set-RowData 139 $row139

Write-Host "Edit complete"
